# Update NATMI Edn1-Ednrb sheet with refreshed TPM-derived values.
# Source data is produced by an external pipeline (no in-sheet formulas),
# so each affected cell is overwritten with its recomputed literal value,
# matching the updated ligand/receptor expression inputs per cluster.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 13.36183833333333
$ws.Cells.Item(2, 8).Value = 40.085515
$ws.Cells.Item(2, 9).Value = 0.8778232942402776
$ws.Cells.Item(2, 10).Value = 0.8933425341492148
$ws.Cells.Item(2, 13).Value = 14.37161333333333
$ws.Cells.Item(2, 14).Value = 43.11484
$ws.Cells.Item(2, 15).Value = 0.4561705932627708
$ws.Cells.Item(2, 16).Value = 0.5019766122855294
$ws.Cells.Item(2, 17).Value = 192.0311739491778
$ws.Cells.Item(2, 18).Value = 1728.2805655426
$ws.Cells.Item(2, 19).Value = 0.4004371729134672
$ws.Cells.Item(2, 20).Value = 0.4484370589027927

# Row 3
$ws.Cells.Item(3, 7).Value = 13.36183833333333
$ws.Cells.Item(3, 8).Value = 40.085515
$ws.Cells.Item(3, 9).Value = 0.8778232942402776
$ws.Cells.Item(3, 10).Value = 0.8933425341492148
$ws.Cells.Item(3, 15).Value = 0.01117178254830525
$ws.Cells.Item(3, 16).Value = 0.01229358849433434
$ws.Cells.Item(3, 17).Value = 4.70291279083
$ws.Cells.Item(3, 18).Value = 42.32621511747
$ws.Cells.Item(3, 19).Value = 0.009806850959089357
$ws.Cells.Item(3, 20).Value = 0.01098238549931627

# Row 4
$ws.Cells.Item(4, 7).Value = 13.36183833333333
$ws.Cells.Item(4, 8).Value = 40.085515
$ws.Cells.Item(4, 9).Value = 0.8778232942402776
$ws.Cells.Item(4, 10).Value = 0.8933425341492148
$ws.Cells.Item(4, 13).Value = 2.798424666666667
$ws.Cells.Item(4, 14).Value = 8.395274000000001
$ws.Cells.Item(4, 15).Value = 0.08882503382091908
$ws.Cells.Item(4, 16).Value = 0.09774433122629669
$ws.Cells.Item(4, 17).Value = 37.39209798401222
$ws.Cells.Item(4, 18).Value = 336.52888185611
$ws.Cells.Item(4, 19).Value = 0.07797268379968327
$ws.Cells.Item(4, 20).Value = 0.08731916855642012

# Row 5
$ws.Cells.Item(5, 7).Value = 13.36183833333333
$ws.Cells.Item(5, 8).Value = 40.085515
$ws.Cells.Item(5, 9).Value = 0.8778232942402776
$ws.Cells.Item(5, 10).Value = 0.8933425341492148
$ws.Cells.Item(5, 13).Value = 8.624592
$ws.Cells.Item(5, 14).Value = 17.249184
$ws.Cells.Item(5, 15).Value = 0.2737539034788959
$ws.Cells.Item(5, 16).Value = 0.2008284606648142
$ws.Cells.Item(5, 17).Value = 115.24040399496
$ws.Cells.Item(5, 18).Value = 691.4424239697601
$ws.Cells.Item(5, 19).Value = 0.2403075533629794
$ws.Cells.Item(5, 20).Value = 0.179408605979591

# Row 6
$ws.Cells.Item(6, 7).Value = 13.36183833333333
$ws.Cells.Item(6, 8).Value = 40.085515
$ws.Cells.Item(6, 9).Value = 0.8778232942402776
$ws.Cells.Item(6, 10).Value = 0.8933425341492148
$ws.Cells.Item(6, 13).Value = 5.358313666666667
$ws.Cells.Item(6, 14).Value = 16.074941
$ws.Cells.Item(6, 15).Value = 0.1700786868891091
$ws.Cells.Item(6, 16).Value = 0.1871570073290255
$ws.Cells.Item(6, 17).Value = 71.59692095329055
$ws.Cells.Item(6, 18).Value = 644.372288579615
$ws.Cells.Item(6, 19).Value = 0.1492990332050585
$ws.Cells.Item(6, 20).Value = 0.1671953152110948

# Row 7
$ws.Cells.Item(7, 9).Value = 0.03473251517428423
$ws.Cells.Item(7, 10).Value = 0.03534655929816115
$ws.Cells.Item(7, 13).Value = 14.37161333333333
$ws.Cells.Item(7, 14).Value = 43.11484
$ws.Cells.Item(7, 15).Value = 0.4561705932627708
$ws.Cells.Item(7, 16).Value = 0.5019766122855294
$ws.Cells.Item(7, 17).Value = 7.598027651906667
$ws.Cells.Item(7, 18).Value = 68.38224886716
$ws.Cells.Item(7, 19).Value = 0.01584395205256143
$ws.Cells.Item(7, 20).Value = 0.01774314609244051

# Row 8
$ws.Cells.Item(8, 9).Value = 0.03473251517428423
$ws.Cells.Item(8, 10).Value = 0.03534655929816115
$ws.Cells.Item(8, 15).Value = 0.01117178254830525
$ws.Cells.Item(8, 16).Value = 0.01229358849433434
$ws.Cells.Item(8, 19).Value = 0.0003880241068828158
$ws.Cells.Item(8, 20).Value = 0.0004345360547021803

# Row 9
$ws.Cells.Item(9, 9).Value = 0.03473251517428423
$ws.Cells.Item(9, 10).Value = 0.03534655929816115
$ws.Cells.Item(9, 13).Value = 2.798424666666667
$ws.Cells.Item(9, 14).Value = 8.395274000000001
$ws.Cells.Item(9, 15).Value = 0.08882503382091908
$ws.Cells.Item(9, 16).Value = 0.09774433122629669
$ws.Cells.Item(9, 17).Value = 1.479479548047334
$ws.Cells.Item(9, 18).Value = 13.315315932426
$ws.Cells.Item(9, 19).Value = 0.003085116835041382
$ws.Cells.Item(9, 20).Value = 0.003454925799749401

# Row 10
$ws.Cells.Item(10, 9).Value = 0.03473251517428423
$ws.Cells.Item(10, 10).Value = 0.03534655929816115
$ws.Cells.Item(10, 13).Value = 8.624592
$ws.Cells.Item(10, 14).Value = 17.249184
$ws.Cells.Item(10, 15).Value = 0.2737539034788959
$ws.Cells.Item(10, 16).Value = 0.2008284606648142
$ws.Cells.Item(10, 17).Value = 4.559675172336
$ws.Cells.Item(10, 18).Value = 27.358051034016
$ws.Cells.Item(10, 19).Value = 0.009508161606600292
$ws.Cells.Item(10, 20).Value = 0.007098595093647278

# Row 11
$ws.Cells.Item(11, 9).Value = 0.03473251517428423
$ws.Cells.Item(11, 10).Value = 0.03534655929816115
$ws.Cells.Item(11, 13).Value = 5.358313666666667
$ws.Cells.Item(11, 14).Value = 16.074941
$ws.Cells.Item(11, 15).Value = 0.1700786868891091
$ws.Cells.Item(11, 16).Value = 0.1871570073290255
$ws.Cells.Item(11, 17).Value = 2.832849344234333
$ws.Cells.Item(11, 18).Value = 25.495644098109
$ws.Cells.Item(11, 19).Value = 0.005907260573198318
$ws.Cells.Item(11, 20).Value = 0.006615356257621779

# Row 12
$ws.Cells.Item(12, 7).Value = 0.3796583333333333
$ws.Cells.Item(12, 8).Value = 1.138975
$ws.Cells.Item(12, 9).Value = 0.02494214647254301
$ws.Cells.Item(12, 10).Value = 0.02538310441646071
$ws.Cells.Item(12, 13).Value = 14.37161333333333
$ws.Cells.Item(12, 14).Value = 43.11484
$ws.Cells.Item(12, 15).Value = 0.4561705932627708
$ws.Cells.Item(12, 16).Value = 0.5019766122855294
$ws.Cells.Item(12, 17).Value = 5.456302765444444
$ws.Cells.Item(12, 18).Value = 49.10672488899999
$ws.Cells.Item(12, 19).Value = 0.01137787375362687
$ws.Cells.Item(12, 20).Value = 0.0127417247642648

# Row 13
$ws.Cells.Item(13, 7).Value = 0.3796583333333333
$ws.Cells.Item(13, 8).Value = 1.138975
$ws.Cells.Item(13, 9).Value = 0.02494214647254301
$ws.Cells.Item(13, 10).Value = 0.02538310441646071
$ws.Cells.Item(13, 15).Value = 0.01117178254830525
$ws.Cells.Item(13, 16).Value = 0.01229358849433434
$ws.Cells.Item(13, 17).Value = 0.13362682495
$ws.Cells.Item(13, 18).Value = 1.20264142455
$ws.Cells.Item(13, 19).Value = 0.0002786482366792293
$ws.Cells.Item(13, 20).Value = 0.0003120494404046885

# Row 14
$ws.Cells.Item(14, 7).Value = 0.3796583333333333
$ws.Cells.Item(14, 8).Value = 1.138975
$ws.Cells.Item(14, 9).Value = 0.02494214647254301
$ws.Cells.Item(14, 10).Value = 0.02538310441646071
$ws.Cells.Item(14, 13).Value = 2.798424666666667
$ws.Cells.Item(14, 14).Value = 8.395274000000001
$ws.Cells.Item(14, 15).Value = 0.08882503382091908
$ws.Cells.Item(14, 16).Value = 0.09774433122629669
$ws.Cells.Item(14, 17).Value = 1.062445244905555
$ws.Cells.Item(14, 18).Value = 9.562007204149999
$ws.Cells.Item(14, 19).Value = 0.002215487003989951
$ws.Cells.Item(14, 20).Value = 0.00248105456563421

# Row 15
$ws.Cells.Item(15, 7).Value = 0.3796583333333333
$ws.Cells.Item(15, 8).Value = 1.138975
$ws.Cells.Item(15, 9).Value = 0.02494214647254301
$ws.Cells.Item(15, 10).Value = 0.02538310441646071
$ws.Cells.Item(15, 13).Value = 8.624592
$ws.Cells.Item(15, 14).Value = 17.249184
$ws.Cells.Item(15, 15).Value = 0.2737539034788959
$ws.Cells.Item(15, 16).Value = 0.2008284606648142
$ws.Cells.Item(15, 17).Value = 3.274398224399999
$ws.Cells.Item(15, 18).Value = 19.6463893464
$ws.Cells.Item(15, 19).Value = 0.006828009958001023
$ws.Cells.Item(15, 20).Value = 0.00509764978685205

# Row 16
$ws.Cells.Item(16, 7).Value = 0.3796583333333333
$ws.Cells.Item(16, 8).Value = 1.138975
$ws.Cells.Item(16, 9).Value = 0.02494214647254301
$ws.Cells.Item(16, 10).Value = 0.02538310441646071
$ws.Cells.Item(16, 13).Value = 5.358313666666667
$ws.Cells.Item(16, 14).Value = 16.074941
$ws.Cells.Item(16, 15).Value = 0.1700786868891091
$ws.Cells.Item(16, 16).Value = 0.1871570073290255
$ws.Cells.Item(16, 17).Value = 2.034328436163888
$ws.Cells.Item(16, 18).Value = 18.308955925475
$ws.Cells.Item(16, 19).Value = 0.00424212752024594
$ws.Cells.Item(16, 20).Value = 0.004750625859304955

# Row 17
$ws.Cells.Item(17, 7).Value = 0.7932915
$ws.Cells.Item(17, 8).Value = 1.586583
$ws.Cells.Item(17, 9).Value = 0.05211631367261799
$ws.Cells.Item(17, 10).Value = 0.03535845997882436
$ws.Cells.Item(17, 13).Value = 14.37161333333333
$ws.Cells.Item(17, 14).Value = 43.11484
$ws.Cells.Item(17, 15).Value = 0.4561705932627708
$ws.Cells.Item(17, 16).Value = 0.5019766122855294
$ws.Cells.Item(17, 17).Value = 11.40087869862
$ws.Cells.Item(17, 18).Value = 68.40527219172
$ws.Cells.Item(17, 19).Value = 0.0237739297267068
$ws.Cells.Item(17, 20).Value = 0.01774911995580373

# Row 18
$ws.Cells.Item(18, 7).Value = 0.7932915
$ws.Cells.Item(18, 8).Value = 1.586583
$ws.Cells.Item(18, 9).Value = 0.05211631367261799
$ws.Cells.Item(18, 10).Value = 0.03535845997882436
$ws.Cells.Item(18, 15).Value = 0.01117178254830525
$ws.Cells.Item(18, 16).Value = 0.01229358849433434
$ws.Cells.Item(18, 17).Value = 0.279211636089
$ws.Cells.Item(18, 18).Value = 1.675269816534
$ws.Cells.Item(18, 19).Value = 0.0005822321235697559
$ws.Cells.Item(18, 20).Value = 0.0004346823567730564

# Row 19
$ws.Cells.Item(19, 7).Value = 0.7932915
$ws.Cells.Item(19, 8).Value = 1.586583
$ws.Cells.Item(19, 9).Value = 0.05211631367261799
$ws.Cells.Item(19, 10).Value = 0.03535845997882436
$ws.Cells.Item(19, 13).Value = 2.798424666666667
$ws.Cells.Item(19, 14).Value = 8.395274000000001
$ws.Cells.Item(19, 15).Value = 0.08882503382091908
$ws.Cells.Item(19, 16).Value = 0.09774433122629669
$ws.Cells.Item(19, 17).Value = 2.219966501457
$ws.Cells.Item(19, 18).Value = 13.319799008742
$ws.Cells.Item(19, 19).Value = 0.004629233324591921
$ws.Cells.Item(19, 20).Value = 0.003456089023821964

# Row 20
$ws.Cells.Item(20, 7).Value = 0.7932915
$ws.Cells.Item(20, 8).Value = 1.586583
$ws.Cells.Item(20, 9).Value = 0.05211631367261799
$ws.Cells.Item(20, 10).Value = 0.03535845997882436
$ws.Cells.Item(20, 13).Value = 8.624592
$ws.Cells.Item(20, 14).Value = 17.249184
$ws.Cells.Item(20, 15).Value = 0.2737539034788959
$ws.Cells.Item(20, 16).Value = 0.2008284606648142
$ws.Cells.Item(20, 17).Value = 6.841815524568
$ws.Cells.Item(20, 18).Value = 27.367262098272
$ws.Cells.Item(20, 19).Value = 0.01426704430280973
$ws.Cells.Item(20, 20).Value = 0.007100985089025734

# Row 21
$ws.Cells.Item(21, 7).Value = 0.7932915
$ws.Cells.Item(21, 8).Value = 1.586583
$ws.Cells.Item(21, 9).Value = 0.05211631367261799
$ws.Cells.Item(21, 10).Value = 0.03535845997882436
$ws.Cells.Item(21, 13).Value = 5.358313666666667
$ws.Cells.Item(21, 14).Value = 16.074941
$ws.Cells.Item(21, 15).Value = 0.1700786868891091
$ws.Cells.Item(21, 16).Value = 0.1871570073290255
$ws.Cells.Item(21, 17).Value = 4.2507046861005
$ws.Cells.Item(21, 18).Value = 25.504228116603
$ws.Cells.Item(21, 19).Value = 0.008863874194939791
$ws.Cells.Item(21, 20).Value = 0.006617583553399885

# Row 22
$ws.Cells.Item(22, 7).Value = 0.158087
$ws.Cells.Item(22, 8).Value = 0.474261
$ws.Cells.Item(22, 9).Value = 0.0103857304402772
$ws.Cells.Item(22, 10).Value = 0.0105693421573389
$ws.Cells.Item(22, 13).Value = 14.37161333333333
$ws.Cells.Item(22, 14).Value = 43.11484
$ws.Cells.Item(22, 15).Value = 0.4561705932627708
$ws.Cells.Item(22, 16).Value = 0.5019766122855294
$ws.Cells.Item(22, 17).Value = 2.271965237026667
$ws.Cells.Item(22, 18).Value = 20.44768713324
$ws.Cells.Item(22, 19).Value = 0.004737664816408468
$ws.Cells.Item(22, 20).Value = 0.005305562570227609

# Row 23
$ws.Cells.Item(23, 7).Value = 0.158087
$ws.Cells.Item(23, 8).Value = 0.474261
$ws.Cells.Item(23, 9).Value = 0.0103857304402772
$ws.Cells.Item(23, 10).Value = 0.0105693421573389
$ws.Cells.Item(23, 15).Value = 0.01117178254830525
$ws.Cells.Item(23, 16).Value = 0.01229358849433434
$ws.Cells.Item(23, 17).Value = 0.055641249042
$ws.Cells.Item(23, 18).Value = 0.500771241378
$ws.Cells.Item(23, 19).Value = 0.0001160271220840914
$ws.Cells.Item(23, 20).Value = 0.0001299351431381443

# Row 24
$ws.Cells.Item(24, 7).Value = 0.158087
$ws.Cells.Item(24, 8).Value = 0.474261
$ws.Cells.Item(24, 9).Value = 0.0103857304402772
$ws.Cells.Item(24, 10).Value = 0.0105693421573389
$ws.Cells.Item(24, 13).Value = 2.798424666666667
$ws.Cells.Item(24, 14).Value = 8.395274000000001
$ws.Cells.Item(24, 15).Value = 0.08882503382091908
$ws.Cells.Item(24, 16).Value = 0.09774433122629669
$ws.Cells.Item(24, 17).Value = 0.4423945602793334
$ws.Cells.Item(24, 18).Value = 3.981551042514
$ws.Cells.Item(24, 19).Value = 0.000922512857612571
$ws.Cells.Item(24, 20).Value = 0.001033093280670995

# Row 25
$ws.Cells.Item(25, 7).Value = 0.158087
$ws.Cells.Item(25, 8).Value = 0.474261
$ws.Cells.Item(25, 9).Value = 0.0103857304402772
$ws.Cells.Item(25, 10).Value = 0.0105693421573389
$ws.Cells.Item(25, 13).Value = 8.624592
$ws.Cells.Item(25, 14).Value = 17.249184
$ws.Cells.Item(25, 15).Value = 0.2737539034788959
$ws.Cells.Item(25, 16).Value = 0.2008284606648142
$ws.Cells.Item(25, 17).Value = 1.363435875504
$ws.Cells.Item(25, 18).Value = 8.180615253024
$ws.Cells.Item(25, 19).Value = 0.002843134248505476
$ws.Cells.Item(25, 20).Value = 0.002122624715698097

# Row 26
$ws.Cells.Item(26, 7).Value = 0.158087
$ws.Cells.Item(26, 8).Value = 0.474261
$ws.Cells.Item(26, 9).Value = 0.0103857304402772
$ws.Cells.Item(26, 10).Value = 0.0105693421573389
$ws.Cells.Item(26, 13).Value = 5.358313666666667
$ws.Cells.Item(26, 14).Value = 16.074941
$ws.Cells.Item(26, 15).Value = 0.1700786868891091
$ws.Cells.Item(26, 16).Value = 0.1871570073290255
$ws.Cells.Item(26, 17).Value = 0.8470797326223334
$ws.Cells.Item(26, 18).Value = 7.623717593601
$ws.Cells.Item(26, 19).Value = 0.001766391395666595
$ws.Cells.Item(26, 20).Value = 0.001978126447604054

